$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1953.3572  # H18 (was 1963.9286)
$ws.Cells.Item(18, 9).Value = 1461.75  # I18 (was 1498.75)
$ws.Cells.Item(18, 11).Value = 1461.75  # K18 (was 1498.75)
$ws.Cells.Item(18, 13).Value = -1177.75  # M18 (was -1214.75)

$ws.Cells.Item(74, 8).Value = 9875  # H74 (was 8250)
$ws.Cells.Item(74, 9).Value = 0  # I74 (was 5000)
$ws.Cells.Item(74, 11).Value = 0  # K74 (was 5000)
$ws.Cells.Item(74, 13).ClearContents()  # M74 (was -4064)

$ws.Cells.Item(77, 8).Value = 9875  # H77 (was 8250)
$ws.Cells.Item(77, 9).Value = 0  # I77 (was 5000)
$ws.Cells.Item(77, 11).Value = 0  # K77 (was 25000)
$ws.Cells.Item(77, 13).ClearContents()  # M77 (was -20320)

$ws.Cells.Item(138, 8).Value = 5799.3286  # H138 (was 5721.5137)
$ws.Cells.Item(138, 9).Value = 3317.2285  # I138 (was 3299.973)
$ws.Cells.Item(138, 11).Value = 9951.6855  # K138 (was 9899.919)
$ws.Cells.Item(138, 13).Value = -4811.6855  # M138 (was -4759.919)

$ws.Cells.Item(141, 8).Value = 2058.111  # H141 (was 2739.2632)
$ws.Cells.Item(141, 10).Value = 0  # J141 (was 15000)
$ws.Cells.Item(141, 12).Value = 0  # L141 (was 45000)
$ws.Cells.Item(141, 14).ClearContents()  # N141 (was -55360)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3553.75  # H2 (was 3361.5715)
$ws.Cells.Item(2, 10).Value = 4899.5  # J2 (was 4899.6665)
$ws.Cells.Item(2, 12).Value = 4899.5  # L2 (was 4899.6665)
$ws.Cells.Item(2, 14).Value = -5125.5  # N2 (was -5125.6665)

$ws.Cells.Item(32, 8).Value = 5621.36  # H32 (was 5532.041)
$ws.Cells.Item(32, 10).Value = 20874.375  # J32 (was 22428.143)
$ws.Cells.Item(32, 12).Value = 20874.375  # L32 (was 22428.143)
$ws.Cells.Item(32, 14).Value = -21448.375  # N32 (was -23002.143)

$ws.Cells.Item(96, 8).Value = 0  # H96 (was 80666.664)
$ws.Cells.Item(96, 10).Value = 0  # J96 (was 80666.664)
$ws.Cells.Item(96, 12).Value = 0  # L96 (was 80666.664)
$ws.Cells.Item(96, 14).ClearContents()  # N96 (was -86158.664)

$ws.Cells.Item(97, 8).Value = 1105.5454  # H97 (was 1002.0769)
$ws.Cells.Item(97, 9).Value = 1116.2  # I97 (was 1002.3333)
$ws.Cells.Item(97, 11).Value = 1116.2  # K97 (was 1002.3333)
$ws.Cells.Item(97, 13).Value = -620.2  # M97 (was -506.3333)

$ws.Cells.Item(116, 8).Value = 3553.75  # H116 (was 3361.5715)
$ws.Cells.Item(116, 10).Value = 4899.5  # J116 (was 4899.6665)
$ws.Cells.Item(116, 12).Value = 4899.5  # L116 (was 4899.6665)
$ws.Cells.Item(116, 14).Value = -9487.5  # N116 (was -9487.666499999999)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3553.75  # H3 (was 3361.5715)
$ws.Cells.Item(3, 10).Value = 4899.5  # J3 (was 4899.6665)
$ws.Cells.Item(3, 12).Value = 4899.5  # L3 (was 4899.6665)
$ws.Cells.Item(3, 14).Value = -5127.5  # N3 (was -5127.6665)

$ws.Cells.Item(92, 8).Value = 0  # H92 (was 25000)
$ws.Cells.Item(92, 10).Value = 0  # J92 (was 25000)
$ws.Cells.Item(92, 12).Value = 0  # L92 (was 25000)
$ws.Cells.Item(92, 14).ClearContents()  # N92 (was -29992)

$ws.Cells.Item(94, 8).Value = 3000  # H94 (was 1750)
$ws.Cells.Item(94, 9).Value = 3000  # I94 (was 1750)
$ws.Cells.Item(94, 11).Value = 3000  # K94 (was 1750)
$ws.Cells.Item(94, 13).Value = -2549  # M94 (was -1299)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4461.5386  # H16 (was 3050.95)
$ws.Cells.Item(16, 9).Value = 4500  # I16 (was 2824.647)
$ws.Cells.Item(16, 11).Value = 4500  # K16 (was 2824.647)
$ws.Cells.Item(16, 13).Value = -4213  # M16 (was -2537.647)

$ws.Cells.Item(21, 8).Value = 4013  # H21 (was 5671)
$ws.Cells.Item(21, 10).Value = 0  # J21 (was 6500)
$ws.Cells.Item(21, 12).Value = 0  # L21 (was 6500)
$ws.Cells.Item(21, 14).ClearContents()  # N21 (was -6970)

$ws.Cells.Item(31, 8).Value = 3322.4092  # H31 (was 4124.107)
$ws.Cells.Item(31, 9).Value = 2694.875  # I31 (was 2933.4285)
$ws.Cells.Item(31, 10).Value = 4995.8335  # J31 (was 5314.7856)
$ws.Cells.Item(31, 11).Value = 2694.875  # K31 (was 2933.4285)
$ws.Cells.Item(31, 12).Value = 4995.8335  # L31 (was 5314.7856)
$ws.Cells.Item(31, 13).Value = -2399.875  # M31 (was -2638.4285)
$ws.Cells.Item(31, 14).Value = -5585.8335  # N31 (was -5904.7856)

$ws.Cells.Item(34, 8).Value = 3322.4092  # H34 (was 4124.107)
$ws.Cells.Item(34, 9).Value = 2694.875  # I34 (was 2933.4285)
$ws.Cells.Item(34, 10).Value = 4995.8335  # J34 (was 5314.7856)
$ws.Cells.Item(34, 11).Value = 2694.875  # K34 (was 2933.4285)
$ws.Cells.Item(34, 12).Value = 4995.8335  # L34 (was 5314.7856)
$ws.Cells.Item(34, 13).Value = -2492.875  # M34 (was -2731.4285)
$ws.Cells.Item(34, 14).Value = -5399.8335  # N34 (was -5718.7856)

$ws.Cells.Item(38, 8).Value = 3036  # H38 (was 3036.3333)
$ws.Cells.Item(38, 9).Value = 3036  # I38 (was 3036.3333)
$ws.Cells.Item(38, 11).Value = 3036  # K38 (was 3036.3333)
$ws.Cells.Item(38, 13).Value = -2659  # M38 (was -2659.3333)

$ws.Cells.Item(46, 8).Value = 3036  # H46 (was 3036.3333)
$ws.Cells.Item(46, 9).Value = 3036  # I46 (was 3036.3333)
$ws.Cells.Item(46, 11).Value = 3036  # K46 (was 3036.3333)
$ws.Cells.Item(46, 13).Value = -2825  # M46 (was -2825.3333)

$ws.Cells.Item(99, 8).Value = 16874.188  # H99 (was 17199.2)
$ws.Cells.Item(99, 9).Value = 14570.143  # I99 (was 14998.667)
$ws.Cells.Item(99, 11).Value = 14570.143  # K99 (was 14998.667)
$ws.Cells.Item(99, 13).Value = -13072.143  # M99 (was -13500.667)

$ws.Cells.Item(107, 8).Value = 2128.3572  # H107 (was 2184.963)
$ws.Cells.Item(107, 10).Value = 1114.6  # J107 (was 1243.25)
$ws.Cells.Item(107, 12).Value = 1114.6  # L107 (was 1243.25)
$ws.Cells.Item(107, 14).Value = -4954.6  # N107 (was -5083.25)

$ws.Cells.Item(113, 8).Value = 4461.5386  # H113 (was 3050.95)
$ws.Cells.Item(113, 9).Value = 4500  # I113 (was 2824.647)
$ws.Cells.Item(113, 11).Value = 4500  # K113 (was 2824.647)
$ws.Cells.Item(113, 13).Value = -2330  # M113 (was -654.6469999999999)

$ws.Cells.Item(126, 8).Value = 16874.188  # H126 (was 17199.2)
$ws.Cells.Item(126, 9).Value = 14570.143  # I126 (was 14998.667)
$ws.Cells.Item(126, 11).Value = 43710.429  # K126 (was 44996.001)
$ws.Cells.Item(126, 13).Value = -41240.429  # M126 (was -42526.001)

$ws.Cells.Item(132, 8).Value = 1578.8823  # H132 (was 1503.381)
$ws.Cells.Item(132, 9).Value = 1344.4615  # I132 (was 1345.2307)
$ws.Cells.Item(132, 10).Value = 2340.75  # J132 (was 1760.375)
$ws.Cells.Item(132, 11).Value = 4033.3845  # K132 (was 4035.6921)
$ws.Cells.Item(132, 12).Value = 7022.25  # L132 (was 5281.125)
$ws.Cells.Item(132, 13).Value = -1503.3845  # M132 (was -1505.6921)
$ws.Cells.Item(132, 14).Value = -12082.25  # N132 (was -10341.125)

$ws.Cells.Item(134, 8).Value = 1915.2572  # H134 (was 1889.8334)
$ws.Cells.Item(134, 9).Value = 1589.2759  # I134 (was 1569.6333)
$ws.Cells.Item(134, 11).Value = 4767.8277  # K134 (was 4708.8999)
$ws.Cells.Item(134, 13).Value = -2232.8277  # M134 (was -2173.8999)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 76962.30499999999  # H2 (was 62528.938)
$ws.Cells.Item(2, 9).Value = 166707.5  # I2 (was 90934.55)
$ws.Cells.Item(2, 10).Value = 37.857143  # J2 (was 36.6)
$ws.Cells.Item(2, 11).Value = 1000245  # K2 (was 545607.3)
$ws.Cells.Item(2, 12).Value = 227.142858  # L2 (was 219.6)
$ws.Cells.Item(2, 13).Value = -1000132  # M2 (was -545494.3)
$ws.Cells.Item(2, 14).Value = -453.142858  # N2 (was -445.6)

$ws.Cells.Item(5, 8).Value = 2645.6667  # H5 (was 2186.6)
$ws.Cells.Item(5, 9).Value = 1470  # I5 (was 1313)
$ws.Cells.Item(5, 10).Value = 4997  # J5 (was 3497)
$ws.Cells.Item(5, 11).Value = 4410  # K5 (was 3939)
$ws.Cells.Item(5, 12).Value = 14991  # L5 (was 10491)
$ws.Cells.Item(5, 13).Value = -4298  # M5 (was -3827)
$ws.Cells.Item(5, 14).Value = -15215  # N5 (was -10715)

$ws.Cells.Item(38, 8).Value = 103.3125  # H38 (was 100.13333)
$ws.Cells.Item(38, 9).Value = 58.88889  # I38 (was 62.666668)
$ws.Cells.Item(38, 10).Value = 160.42857  # J38 (was 156.33333)
$ws.Cells.Item(38, 11).Value = 176.66667  # K38 (was 188.000004)
$ws.Cells.Item(38, 12).Value = 481.28571  # L38 (was 468.99999)
$ws.Cells.Item(38, 13).Value = 170.33333  # M38 (was 158.999996)
$ws.Cells.Item(38, 14).Value = -1175.28571  # N38 (was -1162.99999)

$ws.Cells.Item(135, 8).Value = 2645.6667  # H135 (was 2186.6)
$ws.Cells.Item(135, 9).Value = 1470  # I135 (was 1313)
$ws.Cells.Item(135, 10).Value = 4997  # J135 (was 3497)
$ws.Cells.Item(135, 11).Value = 13230  # K135 (was 11817)
$ws.Cells.Item(135, 12).Value = 44973  # L135 (was 31473)
$ws.Cells.Item(135, 13).Value = -10695  # M135 (was -9282)
$ws.Cells.Item(135, 14).Value = -50043  # N135 (was -36543)

$ws.Cells.Item(139, 8).Value = 1669  # H139 (was 1731.2727)
$ws.Cells.Item(139, 9).Value = 1669  # I139 (was 1704.4)
$ws.Cells.Item(139, 10).Value = 0  # J139 (was 2000)
$ws.Cells.Item(139, 11).Value = 5007  # K139 (was 5113.200000000001)
$ws.Cells.Item(139, 12).Value = 0  # L139 (was 6000)
$ws.Cells.Item(139, 13).Value = 133  # M139 (was 26.79999999999927)
$ws.Cells.Item(139, 14).ClearContents()  # N139 (was -16280)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 12326.167  # H16 (was 12352.429)
$ws.Cells.Item(16, 9).Value = 11991.6  # I16 (was 11993.6)
$ws.Cells.Item(16, 10).Value = 13999  # J16 (was 13249.5)
$ws.Cells.Item(16, 11).Value = 11991.6  # K16 (was 11993.6)
$ws.Cells.Item(16, 12).Value = 13999  # L16 (was 13249.5)
$ws.Cells.Item(16, 13).Value = -11821.6  # M16 (was -11823.6)
$ws.Cells.Item(16, 14).Value = -14339  # N16 (was -13589.5)

$ws.Cells.Item(40, 8).Value = 2649.125  # H40 (was 2685.2222)
$ws.Cells.Item(40, 9).Value = 2741.8572  # I40 (was 2685.2222)
$ws.Cells.Item(40, 10).Value = 2000  # J40 (was 0)
$ws.Cells.Item(40, 11).Value = 2741.8572  # K40 (was 2685.2222)
$ws.Cells.Item(40, 12).Value = 2000  # L40 (was 0)
$ws.Cells.Item(40, 13).Value = -2605.8572  # M40 (was -2549.2222)
$ws.Cells.Item(40, 14).Value = -2272  # N40 (was None)

$ws.Cells.Item(46, 8).Value = 2941.0386  # H46 (was 2970.7585)
$ws.Cells.Item(46, 9).Value = 1785.25  # I46 (was 1947.5)
$ws.Cells.Item(46, 10).Value = 3931.7144  # J46 (was 3925.8)
$ws.Cells.Item(46, 11).Value = 1785.25  # K46 (was 1947.5)
$ws.Cells.Item(46, 12).Value = 3931.7144  # L46 (was 3925.8)
$ws.Cells.Item(46, 13).Value = -1597.25  # M46 (was -1759.5)
$ws.Cells.Item(46, 14).Value = -4307.7144  # N46 (was -4301.8)

$ws.Cells.Item(132, 8).Value = 4449.5557  # H132 (was 4846.5884)
$ws.Cells.Item(132, 10).Value = 4954.875  # J132 (was 5991.2856)
$ws.Cells.Item(132, 12).Value = 14864.625  # L132 (was 17973.8568)
$ws.Cells.Item(132, 14).Value = -19924.625  # N132 (was -23033.8568)

$ws.Cells.Item(134, 8).Value = 30000  # H134 (was 50000)
$ws.Cells.Item(134, 9).Value = 30000  # I134 (was 0)
$ws.Cells.Item(134, 10).Value = 0  # J134 (was 50000)
$ws.Cells.Item(134, 11).Value = 30000  # K134 (was 0)
$ws.Cells.Item(134, 12).Value = 0  # L134 (was 50000)
$ws.Cells.Item(134, 13).Value = -24930  # M134 (was None)
$ws.Cells.Item(134, 14).ClearContents()  # N134 (was -60140)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 3638.75  # H6 (was 3565.5)
$ws.Cells.Item(6, 9).Value = 1516.3334  # I6 (was 2125)
$ws.Cells.Item(6, 10).Value = 10006  # J6 (was 5006)
$ws.Cells.Item(6, 11).Value = 1516.3334  # K6 (was 2125)
$ws.Cells.Item(6, 12).Value = 10006  # L6 (was 5006)
$ws.Cells.Item(6, 13).Value = -1401.3334  # M6 (was -2010)
$ws.Cells.Item(6, 14).Value = -10236  # N6 (was -5236)
